{"js": "// Update the answer cells in the two-digit-division worksheet table.\n// The document has a single table where every 4th row (0, 4, 8, 12, 16)\n// holds five answer cells (\"a\u00f7b=c, d\") and the rows in between are blank\n// spacer rows. We overwrite each answer cell's text in place so the\n// existing run formatting (font / size) on the cell is preserved.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index -> [col0, col1, col2, col3, col4] new values\nconst newValues = {\n  0: [\"27\u00f77=3, 6\", \"83\u00f76=13, 5\", \"44\u00f74=11, 0\", \"23\u00f75=4, 3\", \"65\u00f72=32, 1\"],\n  4: [\"80\u00f77=11, 3\", \"67\u00f76=11, 1\", \"17\u00f72=8, 1\", \"80\u00f75=16, 0\", \"12\u00f74=3, 0\"],\n  8: [\"87\u00f76=14, 3\", \"40\u00f76=6, 4\", \"71\u00f73=23, 2\", \"68\u00f79=7, 5\", \"41\u00f74=10, 1\"],\n  12: [\"10\u00f79=1, 1\", \"18\u00f78=2, 2\", \"73\u00f73=24, 1\", \"63\u00f79=7, 0\", \"67\u00f79=7, 4\"],\n  16: [\"38\u00f73=12, 2\", \"66\u00f72=33, 0\", \"37\u00f75=7, 2\", \"90\u00f77=12, 6\", \"87\u00f72=43, 1\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const rowVals = newValues[rowIndex];\n  for (let col = 0; col < rowVals.length; col++) {\n    const cell = table.getCell(Number(rowIndex), col);\n    cell.value = rowVals[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the answer cells in the two-digit-division worksheet table.\n# The document has a single table where every 4th row (1, 5, 9, 13, 17 in\n# 1-based COM indexing) holds five answer cells (\"a\u00f7b=c, d\") and the rows\n# in between are blank spacer rows. We overwrite each answer cell's\n# Range.Text in place so the existing run formatting (font / size) on the\n# cell is preserved.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"27\u00f77=3, 6\", \"83\u00f76=13, 5\", \"44\u00f74=11, 0\", \"23\u00f75=4, 3\", \"65\u00f72=32, 1\")\n    5  = @(\"80\u00f77=11, 3\", \"67\u00f76=11, 1\", \"17\u00f72=8, 1\", \"80\u00f75=16, 0\", \"12\u00f74=3, 0\")\n    9  = @(\"87\u00f76=14, 3\", \"40\u00f76=6, 4\", \"71\u00f73=23, 2\", \"68\u00f79=7, 5\", \"41\u00f74=10, 1\")\n    13 = @(\"10\u00f79=1, 1\", \"18\u00f78=2, 2\", \"73\u00f73=24, 1\", \"63\u00f79=7, 0\", \"67\u00f79=7, 4\")\n    17 = @(\"38\u00f73=12, 2\", \"66\u00f72=33, 0\", \"37\u00f75=7, 2\", \"90\u00f77=12, 6\", \"87\u00f72=43, 1\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowVals = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowVals.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowVals[$col - 1]\n    }\n}\n"}
